# 13:37 time 15.03.2024 day
# Append 4 new rows (14-17) of registrant data to the sheet, matching the
# source data which stores everything as text (no numeric coercion), and
# grow the used range to A1:G17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while always keeping it as literal
# text (the sheet's existing cells are all plain text -- e.g. "998971300087",
# "3", "32", "15" are stored as text, not numbers). For values that look
# like a number we first flip the cell to the Text number format so the
# assignment isn't auto-coerced into a numeric cell; plain, already
# non-numeric strings don't need that nudge.
function Set-TextValue($Addr, $Value) {
    $looksNumeric = $Value -match '^[0-9]+(\.[0-9]+)?$'

    if ($looksNumeric) {
        $ws.Range($Addr).NumberFormat = "@"
    }
    $ws.Range($Addr).Value = $Value
}

$columns = @("A", "B", "C", "D", "E", "F", "G")

$rows = @(
    @("Ganiyva Nafisaxon Sardor qizi", "998994884859", "None", "994884859", "2", "23 maktab", "8 yosh"),
    @("Isoqov Eldor Fayzullayevich", "998971300087", "eldorisoqov", "971300087", "8", "45", "37"),
    @("Isoqov Eldor Fayzullayevich", "998971300087", "eldorisoqov", "971300087", "8", "45", "12"),
    @("8", "998971300087", "eldorisoqov", "97 130 00 87", "8", "45", "37")
)

$startRow = 14
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $values = $rows[$i]
    for ($j = 0; $j -lt $columns.Length; $j++) {
        $addr = $columns[$j] + $r
        Set-TextValue $addr $values[$j]
    }
}
